$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 6 ("Trend Analysis Node") - tasks below shift up automatically
$ws.Rows.Item(6).Delete()

# Leave selection on the row that took its place, matching the row-delete UX
$ws.Rows.Item(6).Select()

# Re-apply the DueDate formula as one range so Excel groups it as a shared formula
$ws.Range("C2:C8").Formula = '=TEXT(D2,"DD-MMM-YY")'
